$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "27.086.11"
$dCell.ClearFormats()
$ws.Range("E2").Value = "  -0.59%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.826.66"
$dCell.ClearFormats()
$ws.Range("E3").Value = "  +0.25%  "
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.ClearFormats()
$ws.Range("E4").Value = "  -0.35%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "312.81"
$dCell.ClearFormats()
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  -0.33%  "
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4575"
$dCell.ClearFormats()
$ws.Range("E7").Value = "  +6.98%  "
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3728"
$dCell.ClearFormats()
$ws.Range("E8").Value = "  +1.23%  "
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07308"
$dCell.ClearFormats()
$ws.Range("E9").Value = "  +1.01%  "
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.8617"
$dCell.ClearFormats()
$ws.Range("E10").Value = "  -0.07%  "
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "20.90"
$dCell.ClearFormats()
$ws.Range("E11").Value = "  -0.55%  "
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.825.72"
$dCell.ClearFormats()
$ws.Range("E12").Value = "  +0.11%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "6.694"
$dCell.ClearFormats()
$ws.Range("E13").Value = "  +0.26%  "
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "93.08"
$dCell.ClearFormats()
$ws.Range("E14").Value = "  +4.56%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "5.353"
$dCell.ClearFormats()
$ws.Range("E15").Value = "  +0.78%  "
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07103"
$dCell.ClearFormats()
$ws.Range("E16").Value = "  -0.39%  "
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.ClearFormats()
$ws.Range("E17").Value = "  -0.48%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000008848"
$dCell.ClearFormats()
$ws.Range("E18").Value = "  -0.16%  "
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.ClearFormats()
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  -0.35%  "
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "27.118.12"
$dCell.ClearFormats()
$ws.Range("E21").Value = "  -0.59%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.190"
$dCell.ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "10.98"
$dCell.ClearFormats()
$ws.Range("E23").Value = "  +0.89%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "2.001"
$dCell.ClearFormats()
$ws.Range("E24").Value = "  -0.12%  "
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "151.74"
$dCell.ClearFormats()
$ws.Range("E25").Value = "  -1.13%  "
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "2.234"
$dCell.ClearFormats()
$ws.Range("E26").Value = "  +5.04%  "
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "18.52"
$dCell.ClearFormats()
$ws.Range("E27").Value = "  +1.00%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "5.279"
$dCell.ClearFormats()
$ws.Range("E28").Value = "  +0.77%  "
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "117.25"
$dCell.ClearFormats()
$ws.Range("E29").Value = "  +0.82%  "
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08888"
$dCell.ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.195"
$dCell.ClearFormats()
$ws.Range("E31").Value = "  -0.64%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7581"
$dCell.ClearFormats()
$ws.Range("E32").Value = "  -0.29%  "
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "2.962"
$dCell.ClearFormats()
$ws.Range("E33").Value = "  +5.37%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "4.475"
$dCell.ClearFormats()
$ws.Range("E34").Value = "  +0.33%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.ClearFormats()
$ws.Range("E35").Value = "  -0.35%  "
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "1.101"
$dCell.ClearFormats()
$ws.Range("E36").Value = "  -1.23%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "0.01973"
$dCell.ClearFormats()
$ws.Range("E37").Value = "  -0.15%  "
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05268"
$dCell.ClearFormats()
$ws.Range("E38").Value = "  -0.35%  "
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5360"
$dCell.ClearFormats()
$ws.Range("E39").Value = "  +5.91%  "
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "7.192"
$dCell.ClearFormats()
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -0.82%  "
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1718"
$dCell.ClearFormats()
$ws.Range("E42").Value = "  +1.81%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5264"
$dCell.ClearFormats()
$ws.Range("E43").Value = "  +9.58%  "
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "8.586"
$dCell.ClearFormats()
$ws.Range("E44").Value = "  -0.62%  "
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "10.65"
$dCell.ClearFormats()
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +8.83%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "105.81"
$dCell.ClearFormats()
$ws.Range("E47").Value = "  -0.90%  "
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "1.680"
$dCell.ClearFormats()
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  -0.30%  "
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06413"
$dCell.ClearFormats()
$ws.Range("E50").Value = "  +0.18%  "
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "63.52"
$dCell.ClearFormats()
$ws.Range("E51").Value = "  +0.66%  "
